$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing
# username/password columns (A,B) to (B,C).
$ws.Range("A1").EntireColumn.Insert()

# New column A: test-case descriptions
# (values are written in the specific order needed so that new shared
# strings land at the same table indices as the target workbook)
$ws.Range("A1").Value = "Test Case"
$ws.Range("A2").Value = "Valid userId & password"
$ws.Range("A3").Value = "Invalid userId & valid password"
$ws.Range("A5").Value = "Invalid userId & invalid password"
$ws.Range("A4").Value = "Valid userId & invalid password"

# Header row for the shifted columns
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"

# Column width for the new column A (closest achievable value to the
# target raw width of 37.140625 character-units after COM's internal
# pixel-snapping of ColumnWidth)
$ws.Range("A1").EntireColumn.ColumnWidth = 36.33
